$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up two trailing rows (326, 327): the "backup" column (R) goes
# from a blank/placeholder to a completed value of 0 ---
$ws.Cells.Item(326, 18).Value = 0
$ws.Cells.Item(327, 18).Value = 0

# --- An unrelated earlier correction: Q59 (detect_structure) drops back to 0 ---
$ws.Cells.Item(59, 17).Value = 0

# --- Append seven new weekly OHLCV rows (328-334) pulled from the
# upstream stock.yaml job; column R ("backup") is left unset for these,
# matching the not-yet-completed rows' pattern ---
$newRows = @(
    @(45474, 5270,            5624.9501953125,  5253.2998046875,  5552,             5552,             14056044, 2024, 7, 1,  0, 0, 0, 27, 0, 0, 0),
    @(45481, 5580,            5674.75,          5319,             5489.2001953125,  5489.2001953125,  11282380, 2024, 7, 8,  0, 0, 0, 28, 1, 0, 0),
    @(45488, 5495,            5535,             4783.10009765625, 4800.25,          4800.25,          18302228, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 0),
    @(45495, 4799.75,         5073.85009765625, 4510,             4905.39990234375, 4905.39990234375, 19751067, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(45502, 4946.7998046875, 5082,             4631.35009765625, 4695.75,          4695.75,          11799904, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(45509, 4501.0498046875, 4773.5,           4480.10009765625, 4723.89990234375, 4723.89990234375, 10080749, 2024, 8, 5,  0, 0, 0, 32, 0, 0, 0),
    @(45516, 4690,            4811.5,           4593.75,          4769.7998046875,  4769.7998046875,  9301924,  2024, 8, 12, 0, 0, 0, 33, 0, 0, 0)
)

$startRow = 328
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
    # Column A carries the same datetime display format as the rest of the column
    $ws.Range("A$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
